# Update GreatLink Income Bond dividend history with the latest entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DividendHistory")

# Insert 4 new rows right under the header row (row 1), pushing the
# existing history (currently starting at row 2) down to row 6.
$ws.Range("A2:C5").Insert()

# Force the new cells to plain text so dates/numbers aren't
# auto-parsed (matches the rest of the sheet, which stores every
# value - including dates and dividend amounts - as text).
$newRange = $ws.Range("A2:C5")
$newRange.NumberFormat = "@"

# New dividend rows, most recent first.
$newRows = @(
    @("05/03/2025", "05/03/2025", "0.003"),
    @("05/02/2025", "05/02/2025", "0.003"),
    @("06/01/2025", "06/01/2025", "0.003"),
    @("04/12/2024", "04/12/2024", "0.003")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $newRows[$i][2]
}
